$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 596
$ws.Range("F3").Value = 41
$ws.Range("F4").Value = 6400
$ws.Range("F5").Value = 725
$ws.Range("F8").Value = 370
$ws.Range("F11").Value = 701
$ws.Range("F12").Value = 1184
$ws.Range("F13").Value = 82
$ws.Range("F14").Value = 435
$ws.Range("F18").Value = 670
$ws.Range("F21").Value = 82
$ws.Range("F22").Value = 1072
$ws.Range("F24").Value = 2220
$ws.Range("F25").Value = 257
$ws.Range("G25").Value = "不可售"
$ws.Range("F26").Value = 100
$ws.Range("F29").Value = 3585
$ws.Range("F31").Value = 635

$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 183
$ws.Range("F13").Value = 106
$ws.Range("F18").Value = 381
$ws.Range("F20").Value = 4089
$ws.Range("F30").Value = 1

$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 1196
$ws.Range("F6").Value = 1576
$ws.Range("F7").Value = 432
$ws.Range("F10").Value = 774

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 1196
$ws.Range("F4").Value = 1576
$ws.Range("F5").Value = 432
$ws.Range("F7").Value = 774
$ws.Range("F8").Value = 596
$ws.Range("F9").Value = 183
$ws.Range("F10").Value = 6400
$ws.Range("F12").Value = 725
$ws.Range("F16").Value = 370
$ws.Range("F19").Value = 701
$ws.Range("F21").Value = 106
$ws.Range("F22").Value = 106
$ws.Range("F24").Value = 1184
$ws.Range("F25").Value = 82
$ws.Range("F27").Value = 381
$ws.Range("F32").Value = 670
$ws.Range("F37").Value = 1072
$ws.Range("F39").Value = 2220
$ws.Range("B40").Value = "'2024-04-20"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "上海·冰兔2024线下live《过去和未来》"
$ws.Range("D40").Value = "重庆南路308号3楼 上海MAO LIVEHOUSE"
$ws.Range("E40").Value = "2024.04.20 13:00-04.20 15:00"
$ws.Range("F40").Value = 207
$ws.Range("G40").Value = 198
$ws.Range("H40").Value = "https://show.bilibili.com/platform/detail.html?id=81654"
$ws.Range("I40").Value = "//i1.hdslb.com/bfs/openplatform/202402/OEHnMZmi1706851347869.jpeg"
$ws.Range("B41").Value = "'2024-04-21"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "上海·今泉爱夏  巡演"
$ws.Range("D41").Value = "瑞虹路188号3楼 Modernsky Lab"
$ws.Range("E41").Value = "2024.04.21 20:00-04.21 21:30"
$ws.Range("F41").Value = 33
$ws.Range("G41").Value = 328
$ws.Range("H41").Value = "https://show.bilibili.com/platform/detail.html?id=81891"
$ws.Range("I41").Value = "//i1.hdslb.com/bfs/openplatform/202402/YJENeaUi1708313389899.jpeg"
$ws.Range("C42").Value = "上海· 夏川里美 2024 巡回演唱会 出道 25 周年纪念专场"
$ws.Range("D42").Value = "东大名路889号 友邦大剧院"
$ws.Range("E42").Value = "2024.04.26 19:30-04.26 21:30"
$ws.Range("F42").Value = 47
$ws.Range("H42").Value = "https://show.bilibili.com/platform/detail.html?id=81139"
$ws.Range("I42").Value = "//i2.hdslb.com/bfs/openplatform/202401/0Fj4cYOH1705652393930.jpeg"
$ws.Range("B44").Value = "'2024-04-26"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "上海·「angela LIVE 2024」in SHANGHAI"
$ws.Range("D44").Value = "万航渡后路19号3楼 瓦肆VAS SHANGHAI"
$ws.Range("E44").Value = "2024.04.26 19:00-04.26 20:30"
$ws.Range("F44").Value = 1661
$ws.Range("G44").Value = 480
$ws.Range("H44").Value = "https://show.bilibili.com/platform/detail.html?id=82039"
$ws.Range("I44").Value = "//i2.hdslb.com/bfs/openplatform/202402/H9L22d9R1708678603570.jpeg"
$ws.Range("F45").Value = 100
$ws.Range("F47").Value = 3585
$ws.Range("F51").Value = 635
